$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "25.772.06"
$ws.Range("E2").Value = "  -0.19%  "

# Row 3
$ws.Range("D3").Value = "1.635.08"
$ws.Range("E3").Value = "  +0.03%  "

# Row 4
$ws.Range("E4").Value = "  -0.17%  "

# Row 5
$ws.Range("D5").Formula = "'215.46"
$ws.Range("E5").Value = "  +0.16%  "

# Row 6
$ws.Range("E6").Value = "  -0.69%  "

# Row 7
$ws.Range("E7").Value = "  -0.16%  "

# Row 8
$ws.Range("E8").Value = "  -0.29%  "

# Row 9
$ws.Range("E9").Value = "  -1.16%  "

# Row 10
$ws.Range("D10").Formula = "'19.55"
$ws.Range("E10").Value = "  -1.83%  "

# Row 11
$ws.Range("D11").Formula = "'0.0792"
$ws.Range("E11").Value = "  +1.58%  "

# Row 12
$ws.Range("E12").Value = "  +0.16%  "

# Row 13
$ws.Range("D13").Value = "1.860.88"
$ws.Range("E13").Value = "  +0.05%  "

# Row 14
$ws.Range("D14").Value = "1.641.28"
$ws.Range("E14").Value = "  +0.39%  "

# Row 15
$ws.Range("E15").Value = "  +0.72%  "

# Row 16
$ws.Range("E16").Value = "  -0.65%  "

# Row 17
$ws.Range("D17").Formula = "'63.19"
$ws.Range("E17").Value = "  +0.28%  "

# Row 18
$ws.Range("D18").Value = "25.805.07"
$ws.Range("E18").Value = "  -0.08%  "

# Row 19
$ws.Range("E19").Value = "  -0.18%  "

# Row 21
$ws.Range("D21").Formula = "'192.42"
$ws.Range("E21").Value = "  -0.64%  "

# Row 22
$ws.Range("E22").Value = "  +0.60%  "

# Row 23
$ws.Range("D23").Formula = "'6.29"
$ws.Range("E23").Value = "  +2.26%  "

# Row 24
$ws.Range("E24").Value = "  +4.49%  "

# Row 25
$ws.Range("E25").Value = "  -0.11%  "

# Row 26
$ws.Range("D26").Formula = "'141.48"
$ws.Range("E26").Value = "  +1.29%  "

# Row 27
$ws.Range("E27").Value = "  +1.34%  "

# Row 28
$ws.Range("E28").Value = "  +1.10%  "

# Row 29
$ws.Range("E29").Value = "  +0.35%  "

# Row 30
$ws.Range("E30").Value = "  -0.14%  "

# Row 31
$ws.Range("D31").Formula = "'0.0493"
$ws.Range("E31").Value = "  -0.35%  "

# Row 32
$ws.Range("E32").Value = "  -0.06%  "

# Row 33
$ws.Range("E33").Value = "  -0.53%  "

# Row 34
$ws.Range("E34").Value = "  -0.46%  "

# Row 35
$ws.Range("E35").Value = "  -0.29%  "

# Row 36
$ws.Range("D36").Formula = "'0.902"
$ws.Range("E36").Value = "  +0.34%  "

# Row 37
$ws.Range("D37").Value = "1.131.72"
$ws.Range("E37").Value = "  +1.40%  "

# Row 38
$ws.Range("E38").Value = "  -2.08%  "

# Row 39
$ws.Range("D39").Formula = "'0.543"
$ws.Range("E39").Value = "  -0.85%  "

# Row 40
$ws.Range("E40").Value = "  -0.75%  "

# Row 41
$ws.Range("E41").Value = "  +0.08%  "

# Row 42
$ws.Range("E42").Value = "  +0.59%  "

# Row 43
$ws.Range("E43").Value = "  +0.50%  "

# Row 44
$ws.Range("D44").Formula = "'100.72"
$ws.Range("E44").Value = "  +1.36%  "

# Row 45
$ws.Range("D45").Formula = "'0.804"
$ws.Range("E45").Value = "  +0.58%  "

# Row 46
$ws.Range("D46").Value = "1.769.98"
$ws.Range("E46").Value = "  -0.07%  "

# Row 47
$ws.Range("B47").Value = "Aave"
$ws.Range("C47").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D47").Formula = "'55.31"
$ws.Range("E47").Value = "  -0.35%  "

# Row 48
$ws.Range("B48").Value = "Mantle"
$ws.Range("C48").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D48").Formula = "'0.417"
$ws.Range("E48").Value = "  -0.95%  "

# Row 49
$ws.Range("B49").Value = "RenderToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D49").Formula = "'1.44"
$ws.Range("E49").Value = "  +4.80%  "

# Row 50
$ws.Range("E50").Value = "  -0.16%  "

# Row 51
$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").Formula = "'7.47"
$ws.Range("E51").Value = "  -2.30%  "
